$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row -> new F value }
$updates = @{
    "展览" = @{
        2  = 65
        3  = 366
        5  = 246
        6  = 13578
        9  = 5499
        10 = 567
        15 = 62
        17 = 730
        19 = 9153
        21 = 3683
    }
    "全部类型" = @{
        2  = 65
        3  = 366
        6  = 246
        7  = 13578
        10 = 5499
        11 = 567
        16 = 62
        18 = 730
        21 = 9153
        23 = 3683
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $updates[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $ws.Range("F$row").Value = $rowMap[$row]
    }
}
